$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1534.8064
$ws.Range("I98").Value = 1295.6552
$ws.Range("K98").Value = 1295.6552
$ws.Range("M98").Value = 202.3448000000001
$ws.Range("H116").Value = 10284.444
$ws.Range("J116").Value = 10527.2
$ws.Range("L116").Value = 10527.2
$ws.Range("N116").Value = -17411.2
$ws.Range("H122").Value = 1534.8064
$ws.Range("I122").Value = 1295.6552
$ws.Range("K122").Value = 3886.9656
$ws.Range("M122").Value = -1436.9656
$ws.Range("H125").Value = 2967.5
$ws.Range("J125").Value = 2967.5
$ws.Range("L125").Value = 26707.5
$ws.Range("N125").Value = -31627.5
$ws.Range("H137").Value = 10002244
$ws.Range("I137").Value = 14287800
$ws.Range("J137").Value = 2613.6
$ws.Range("K137").Value = 42863400
$ws.Range("L137").Value = 7840.799999999999
$ws.Range("M137").Value = -42860850
$ws.Range("N137").Value = -12940.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5253.0425
$ws.Range("I32").Value = 5127.6343
$ws.Range("K32").Value = 5127.6343
$ws.Range("M32").Value = -4840.6343
$ws.Range("H45").Value = 3462.5293
$ws.Range("I45").Value = 2375
$ws.Range("J45").Value = 4055.7273
$ws.Range("K45").Value = 2375
$ws.Range("L45").Value = 4055.7273
$ws.Range("M45").Value = -1998
$ws.Range("N45").Value = -4809.7273
$ws.Range("H61").Value = 30871290
$ws.Range("I61").Value = 33335156
$ws.Range("K61").Value = 33335156
$ws.Range("M61").Value = -33334944
$ws.Range("H63").Value = 979.5
$ws.Range("I63").Value = 979.5
$ws.Range("K63").Value = 979.5
$ws.Range("M63").Value = -293.5
$ws.Range("H66").Value = 979.5
$ws.Range("I66").Value = 979.5
$ws.Range("K66").Value = 4897.5
$ws.Range("M66").Value = -1465.5
$ws.Range("H74").Value = 1725.4117
$ws.Range("I74").Value = 1420.3928
$ws.Range("K74").Value = 1420.3928
$ws.Range("M74").Value = -546.3928000000001
$ws.Range("H77").Value = 1725.4117
$ws.Range("I77").Value = 1420.3928
$ws.Range("K77").Value = 7101.964
$ws.Range("M77").Value = -2733.964
$ws.Range("H97").Value = 1103.1212
$ws.Range("I97").Value = 1014.1429
$ws.Range("J97").Value = 1258.8334
$ws.Range("K97").Value = 1014.1429
$ws.Range("L97").Value = 1258.8334
$ws.Range("M97").Value = -518.1429000000001
$ws.Range("N97").Value = -2250.8334
$ws.Range("M122").Value = -4163.674
$ws.Range("H122").Value = 2204.558
$ws.Range("I122").Value = 2204.558
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6613.674
$ws.Range("L122").Value = 0
$ws.Range("H132").Value = 13893130
$ws.Range("I132").Value = 4260.25
$ws.Range("K132").Value = 12780.75
$ws.Range("M132").Value = -10250.75
$ws.Range("H136").Value = 30871290
$ws.Range("I136").Value = 33335156
$ws.Range("K136").Value = 100005468
$ws.Range("M136").Value = -100002918
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1170335
$ws.Range("I86").Value = 1945676.9
$ws.Range("K86").Value = 1945676.9
$ws.Range("M86").Value = -1944553.9
$ws.Range("H89").Value = 1170335
$ws.Range("I89").Value = 1945676.9
$ws.Range("K89").Value = 9728384.5
$ws.Range("M89").Value = -9722768.5
$ws.Range("H105").Value = 446799.72
$ws.Range("I105").Value = 758848.4399999999
$ws.Range("J105").Value = 4730.75
$ws.Range("K105").Value = 758848.4399999999
$ws.Range("L105").Value = 4730.75
$ws.Range("M105").Value = -757101.4399999999
$ws.Range("N105").Value = -8224.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 7541
$ws.Range("I62").Value = 7832.6665
$ws.Range("K62").Value = 7832.6665
$ws.Range("M62").Value = -7208.6665
$ws.Range("H65").Value = 7541
$ws.Range("I65").Value = 7832.6665
$ws.Range("K65").Value = 39163.3325
$ws.Range("M65").Value = -36043.3325
$ws.Range("H132").Value = 1996.0555
$ws.Range("I132").Value = 1967.4375
$ws.Range("K132").Value = 5902.3125
$ws.Range("M132").Value = -3372.3125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1447.6
$ws.Range("J97").Value = 1219.7778
$ws.Range("L97").Value = 3659.3334
$ws.Range("N97").Value = -4651.3334
$ws.Range("H134").Value = 13053100
$ws.Range("I134").Value = 16669777
$ws.Range("K134").Value = 50009331
$ws.Range("M134").Value = -50004261
$ws.Range("H139").Value = 1966.6813
$ws.Range("I139").Value = 2129.8823
$ws.Range("K139").Value = 6389.646900000001
$ws.Range("M139").Value = -1249.646900000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2136.1875
$ws.Range("I122").Value = 2045.1538
$ws.Range("J122").Value = 2530.6667
$ws.Range("K122").Value = 6135.4614
$ws.Range("L122").Value = 7592.000100000001
$ws.Range("M122").Value = -3685.4614
$ws.Range("N122").Value = -12492.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 4214.857
$ws.Range("I82").Value = 2992.6
$ws.Range("J82").Value = 5326
$ws.Range("K82").Value = 2992.6
$ws.Range("L82").Value = 5326
$ws.Range("M82").Value = -2631.6
$ws.Range("N82").Value = -6048
$ws.Range("H85").Value = 4214.857
$ws.Range("I85").Value = 2992.6
$ws.Range("J85").Value = 5326
$ws.Range("K85").Value = 2992.6
$ws.Range("L85").Value = 5326
$ws.Range("M85").Value = -1744.6
$ws.Range("N85").Value = -7822
$ws.Range("H93").Value = 1504177.4
$ws.Range("I93").Value = 942.88
$ws.Range("K93").Value = 942.88
$ws.Range("M93").Value = 305.12
$ws.Range("H100").Value = 15644664
$ws.Range("I100").Value = 4021.125
$ws.Range("K100").Value = 4021.125
$ws.Range("M100").Value = -3480.125
$ws.Range("H122").Value = 3025.2432
$ws.Range("I122").Value = 2796.5715
$ws.Range("J122").Value = 7027
$ws.Range("K122").Value = 8389.7145
$ws.Range("L122").Value = 21081
$ws.Range("M122").Value = -5939.7145
$ws.Range("N122").Value = -25981
$ws.Range("H136").Value = 2623.16
$ws.Range("I136").Value = 2397.5293
$ws.Range("J136").Value = 3102.625
$ws.Range("K136").Value = 7192.5879
$ws.Range("L136").Value = 9307.875
$ws.Range("M136").Value = -4642.5879
$ws.Range("N136").Value = -14407.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 12218.143
$ws.Range("I45").Value = 24984.5
$ws.Range("J45").Value = 7111.6
$ws.Range("K45").Value = 24984.5
$ws.Range("L45").Value = 7111.6
$ws.Range("M45").Value = -24493.5
$ws.Range("N45").Value = -8093.6
$ws.Range("H54").Value = 42000
$ws.Range("J54").Value = 42000
$ws.Range("L54").Value = 42000
$ws.Range("N54").Value = -43040
$ws.Range("M81").Value = -1346
$ws.Range("H81").Value = 1203.5
$ws.Range("I81").Value = 1203.5
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2407
$ws.Range("L81").Value = 0
$ws.Range("M84").Value = -6731
$ws.Range("H84").Value = 1203.5
$ws.Range("I84").Value = 1203.5
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 12035
$ws.Range("L84").Value = 0
$ws.Range("H86").Value = 78999
$ws.Range("J86").Value = 78999
$ws.Range("L86").Value = 78999
$ws.Range("N86").Value = -81245
$ws.Range("H89").Value = 78999
$ws.Range("J89").Value = 78999
$ws.Range("L89").Value = 394995
$ws.Range("N89").Value = -406227
$ws.Range("N81").ClearContents()
$ws.Range("N84").ClearContents()
